# Apply the updated "location01" result values to column D of Sheet1.
# These values reflect a fix to the NSGA multiobjective example's tilting
# calculation, as described in the commit message.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    "D2"  = 3090.974289607632
    "D3"  = 3090.97428960763
    "D4"  = 3090.97428960763
    "D6"  = 8831.355113164813
    "D7"  = 8831.355113164813
    "D19" = 79354.61403309148
    "D20" = 79354.61403309148
    "D24" = 149626.5932405632
    "D25" = 149626.593240563
    "D28" = 2992.531864811288
    "D29" = 2992.531864811288
    "D30" = 149626.593240563
    "D38" = -3390.728536028876
    "D39" = -3390.728536028873
    "D41" = 3390.728536028873
    "D42" = 67814.57072057677
    "D43" = 67814.57072057677
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
